$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These columns (runs/balls/fours/sixes) store numeric-looking values as
# text in the source workbook; prefix with an apostrophe so Excel keeps
# them as text instead of coercing to numeric cells.

# Row 2: runs/balls/fours updated
$ws.Range("C2").Value = "'13"
$ws.Range("D2").Value = "'8"
$ws.Range("E2").Value = "'2"

# Row 3: runs/balls updated
$ws.Range("C3").Value = "'4"
$ws.Range("D3").Value = "'5"

# Row 4: runs/balls/fours updated
$ws.Range("C4").Value = "'7"
$ws.Range("D4").Value = "'12"
$ws.Range("E4").Value = "'0"

# Row 6: runs/balls updated
$ws.Range("C6").Value = "'5"
$ws.Range("D6").Value = "'6"
